$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Il12a"
$ws.Cells.Item(2,3).Value = "Il12rb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.019271
$ws.Cells.Item(2,8).Value = 0.057813
$ws.Cells.Item(2,9).Value = 0.009412517461340836
$ws.Cells.Item(2,10).Value = 0.009412517461340836
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.3754836666666666
$ws.Cells.Item(2,14).Value = 1.126451
$ws.Cells.Item(2,15).Value = 0.1282941949732911
$ws.Cells.Item(2,16).Value = 0.1282941949732911
$ws.Cells.Item(2,17).Value = 0.007235945740333333
$ws.Cells.Item(2,18).Value = 0.065123511663
$ws.Cells.Item(2,19).Value = 0.001207571350374768
$ws.Cells.Item(2,20).Value = 0.001207571350374768

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Il12a"
$ws.Cells.Item(3,3).Value = "Il12rb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.019271
$ws.Cells.Item(3,8).Value = 0.057813
$ws.Cells.Item(3,9).Value = 0.009412517461340836
$ws.Cells.Item(3,10).Value = 0.009412517461340836
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.544463
$ws.Cells.Item(3,14).Value = 1.633389
$ws.Cells.Item(3,15).Value = 0.1860305746394907
$ws.Cells.Item(3,16).Value = 0.1860305746394907
$ws.Cells.Item(3,17).Value = 0.010492346473
$ws.Cells.Item(3,18).Value = 0.094431118257
$ws.Cells.Item(3,19).Value = 0.001751016032137476
$ws.Cells.Item(3,20).Value = 0.001751016032137476

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Il12a"
$ws.Cells.Item(4,3).Value = "Il12rb2"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.019271
$ws.Cells.Item(4,8).Value = 0.057813
$ws.Cells.Item(4,9).Value = 0.009412517461340836
$ws.Cells.Item(4,10).Value = 0.009412517461340836
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.9905483333333333
$ws.Cells.Item(4,14).Value = 2.971645
$ws.Cells.Item(4,15).Value = 0.3384477469693805
$ws.Cells.Item(4,16).Value = 0.3384477469693805
$ws.Cells.Item(4,17).Value = 0.01908885693166667
$ws.Cells.Item(4,18).Value = 0.171799712385
$ws.Cells.Item(4,19).Value = 0.003185645328100758
$ws.Cells.Item(4,20).Value = 0.003185645328100758

$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Il12a"
$ws.Cells.Item(5,3).Value = "Il12rb2"
$ws.Cells.Item(5,4).Value = "Resolving-Mac"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.019271
$ws.Cells.Item(5,8).Value = 0.057813
$ws.Cells.Item(5,9).Value = 0.009412517461340836
$ws.Cells.Item(5,10).Value = 0.009412517461340836
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.016244333333333
$ws.Cells.Item(5,14).Value = 3.048733
$ws.Cells.Item(5,15).Value = 0.3472274834178377
$ws.Cells.Item(5,16).Value = 0.3472274834178377
$ws.Cells.Item(5,17).Value = 0.01958404454766667
$ws.Cells.Item(5,18).Value = 0.176256400929
$ws.Cells.Item(5,19).Value = 0.003268284750727833
$ws.Cells.Item(5,20).Value = 0.003268284750727833

$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,2).Value = "Il12a"
$ws.Cells.Item(6,3).Value = "Il12rb2"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.028109
$ws.Cells.Item(6,8).Value = 6.084327
$ws.Cells.Item(6,9).Value = 0.9905874825386592
$ws.Cells.Item(6,10).Value = 0.9905874825386591
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.3754836666666666
$ws.Cells.Item(6,14).Value = 1.126451
$ws.Cells.Item(6,15).Value = 0.1282941949732911
$ws.Cells.Item(6,16).Value = 0.1282941949732911
$ws.Cells.Item(6,17).Value = 0.7615218037196666
$ws.Cells.Item(6,18).Value = 6.853696233476999
$ws.Cells.Item(6,19).Value = 0.1270866236229163
$ws.Cells.Item(6,20).Value = 0.1270866236229163

$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Il12a"
$ws.Cells.Item(7,3).Value = "Il12rb2"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.028109
$ws.Cells.Item(7,8).Value = 6.084327
$ws.Cells.Item(7,9).Value = 0.9905874825386592
$ws.Cells.Item(7,10).Value = 0.9905874825386591
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.544463
$ws.Cells.Item(7,14).Value = 1.633389
$ws.Cells.Item(7,15).Value = 0.1860305746394907
$ws.Cells.Item(7,16).Value = 0.1860305746394907
$ws.Cells.Item(7,17).Value = 1.104230310467
$ws.Cells.Item(7,18).Value = 9.938072794203
$ws.Cells.Item(7,19).Value = 0.1842795586073533
$ws.Cells.Item(7,20).Value = 0.1842795586073533

$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Il12a"
$ws.Cells.Item(8,3).Value = "Il12rb2"
$ws.Cells.Item(8,4).Value = "MuSCs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.028109
$ws.Cells.Item(8,8).Value = 6.084327
$ws.Cells.Item(8,9).Value = 0.9905874825386592
$ws.Cells.Item(8,10).Value = 0.9905874825386591
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.9905483333333333
$ws.Cells.Item(8,14).Value = 2.971645
$ws.Cells.Item(8,15).Value = 0.3384477469693805
$ws.Cells.Item(8,16).Value = 0.3384477469693805
$ws.Cells.Item(8,17).Value = 2.008939989768333
$ws.Cells.Item(8,18).Value = 18.080459907915
$ws.Cells.Item(8,19).Value = 0.3352621016412797
$ws.Cells.Item(8,20).Value = 0.3352621016412797

$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Il12a"
$ws.Cells.Item(9,3).Value = "Il12rb2"
$ws.Cells.Item(9,4).Value = "Resolving-Mac"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.028109
$ws.Cells.Item(9,8).Value = 6.084327
$ws.Cells.Item(9,9).Value = 0.9905874825386592
$ws.Cells.Item(9,10).Value = 0.9905874825386591
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.016244333333333
$ws.Cells.Item(9,14).Value = 3.048733
$ws.Cells.Item(9,15).Value = 0.3472274834178377
$ws.Cells.Item(9,16).Value = 0.3472274834178377
$ws.Cells.Item(9,17).Value = 2.061054278632334
$ws.Cells.Item(9,18).Value = 18.549488507691
$ws.Cells.Item(9,19).Value = 0.3439591986671098
$ws.Cells.Item(9,20).Value = 0.3439591986671099
